# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview zh-cn/de-de status columns, and the Status column on
#   each per-locale report sheet).
# - Narrow the "Latest Handoff Datetime"-style status columns (Overview!E:F,
#   and the Status column C on each locale sheet) to their new width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Excel's ColumnWidth (character units) is quantized internally to 1/6ths, so
# the nearest achievable setting for the target stored width (13.4101845877511)
# is 12.5, which rounds to a stored width of 13.333333333333334.
$newColumnWidth = 12.5

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns("E:F").ColumnWidth = $newColumnWidth

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns("C:C").ColumnWidth = $newColumnWidth

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns("C:C").ColumnWidth = $newColumnWidth
